$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$updates = @(
    @{ Cell = 'D2'; Value = '28.709.74' }
    @{ Cell = 'E2'; Value = '  -1.62%  ' }
    @{ Cell = 'D3'; Value = '1.802.74' }
    @{ Cell = 'E3'; Value = '  -1.29%  ' }
    @{ Cell = 'E4'; Value = '  +0.33%  ' }
    @{ Cell = 'D5'; Value = '231.72' }
    @{ Cell = 'E5'; Value = '  -2.08%  ' }
    @{ Cell = 'D6'; Value = '0.5948' }
    @{ Cell = 'E6'; Value = '  -2.87%  ' }
    @{ Cell = 'D7'; Value = '1.004' }
    @{ Cell = 'E7'; Value = '  +0.36%  ' }
    @{ Cell = 'D8'; Value = '0.2777' }
    @{ Cell = 'E8'; Value = '  -1.41%  ' }
    @{ Cell = 'D9'; Value = '0.06830' }
    @{ Cell = 'E9'; Value = '  -3.86%  ' }
    @{ Cell = 'D10'; Value = '23.33' }
    @{ Cell = 'E10'; Value = '  -0.88%  ' }
    @{ Cell = 'D12'; Value = '1.798.13' }
    @{ Cell = 'E12'; Value = '  -1.55%  ' }
    @{ Cell = 'D13'; Value = '4.791' }
    @{ Cell = 'E13'; Value = '  -0.55%  ' }
    @{ Cell = 'D14'; Value = '0.6246' }
    @{ Cell = 'E14'; Value = '  -1.20%  ' }
    @{ Cell = 'D15'; Value = '2.048.74' }
    @{ Cell = 'E15'; Value = '  -1.20%  ' }
    @{ Cell = 'D16'; Value = '0.000009273' }
    @{ Cell = 'E16'; Value = '  -7.95%  ' }
    @{ Cell = 'D17'; Value = '75.35' }
    @{ Cell = 'E17'; Value = '  -4.57%  ' }
    @{ Cell = 'D18'; Value = '28.688.88' }
    @{ Cell = 'E18'; Value = '  -1.68%  ' }
    @{ Cell = 'D19'; Value = '5.470' }
    @{ Cell = 'E19'; Value = '  -6.81%  ' }
    @{ Cell = 'E20'; Value = '  +0.33%  ' }
    @{ Cell = 'D21'; Value = '210.28' }
    @{ Cell = 'D22'; Value = '11.45' }
    @{ Cell = 'E22'; Value = '  -2.83%  ' }
    @{ Cell = 'D23'; Value = '6.858' }
    @{ Cell = 'E23'; Value = '  -2.10%  ' }
    @{ Cell = 'D24'; Value = '1.004' }
    @{ Cell = 'E24'; Value = '  +0.30%  ' }
    @{ Cell = 'D25'; Value = '154.27' }
    @{ Cell = 'E25'; Value = '  -0.47%  ' }
    @{ Cell = 'D26'; Value = '7.845' }
    @{ Cell = 'E26'; Value = '  -2.41%  ' }
    @{ Cell = 'D27'; Value = '0.1275' }
    @{ Cell = 'E27'; Value = '  -3.19%  ' }
    @{ Cell = 'D28'; Value = '16.40' }
    @{ Cell = 'E28'; Value = '  -1.27%  ' }
    @{ Cell = 'D29'; Value = '1.432' }
    @{ Cell = 'E29'; Value = '  -3.80%  ' }
    @{ Cell = 'D30'; Value = '0.06172' }
    @{ Cell = 'E30'; Value = '  -2.69%  ' }
    @{ Cell = 'E31'; Value = '  -2.24%  ' }
    @{ Cell = 'D32'; Value = '3.782' }
    @{ Cell = 'E32'; Value = '  -1.04%  ' }
    @{ Cell = 'E33'; Value = '  -1.21%  ' }
    @{ Cell = 'E34'; Value = '  -1.72%  ' }
    @{ Cell = 'D35'; Value = '1.062' }
    @{ Cell = 'E35'; Value = '  -5.77%  ' }
    @{ Cell = 'D36'; Value = '0.6404' }
    @{ Cell = 'E36'; Value = '  -1.40%  ' }
    @{ Cell = 'D37'; Value = '2.496' }
    @{ Cell = 'E37'; Value = '  -1.86%  ' }
    @{ Cell = 'D38'; Value = '2.715' }
    @{ Cell = 'E38'; Value = '  -1.19%  ' }
    @{ Cell = 'D39'; Value = '0.01711' }
    @{ Cell = 'E39'; Value = '  -1.80%  ' }
    @{ Cell = 'E40'; Value = '  -2.35%  ' }
    @{ Cell = 'D41'; Value = '1.132.96' }
    @{ Cell = 'E41'; Value = '  -6.87%  ' }
    @{ Cell = 'D42'; Value = '0.8676' }
    @{ Cell = 'E42'; Value = '  -6.20%  ' }
    @{ Cell = 'D43'; Value = '1.004' }
    @{ Cell = 'E43'; Value = '  +0.47%  ' }
    @{ Cell = 'E44'; Value = '  -0.64%  ' }
    @{ Cell = 'D45'; Value = '1.966.42' }
    @{ Cell = 'E45'; Value = '  -0.36%  ' }
    @{ Cell = 'D46'; Value = '60.60' }
    @{ Cell = 'E46'; Value = '  -3.83%  ' }
    @{ Cell = 'D47'; Value = '0.00000000112' }
    @{ Cell = 'E47'; Value = '  -5.83%  ' }
    @{ Cell = 'D48'; Value = '1.598' }
    @{ Cell = 'E48'; Value = '  -1.58%  ' }
    @{ Cell = 'B49'; Value = 'EnergySwap' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' }
    @{ Cell = 'D49'; Value = '8.356' }
    @{ Cell = 'E49'; Value = '  -3.17%  ' }
    @{ Cell = 'B50'; Value = 'Cronos' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' }
    @{ Cell = 'D50'; Value = '0.05471' }
    @{ Cell = 'E50'; Value = '  -0.86%  ' }
    @{ Cell = 'D51'; Value = '0.4494' }
    @{ Cell = 'E51'; Value = '  -1.50%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.Value = "'" + $u.Value
    $cell.Style = "Normal"
}
